$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Berenjena" (eggplant) series
# for "Terminal La Palmera de La Serena". It belongs right after the
# existing row 317 (chronologically between the 44979 and 45222 dates),
# so insert a fresh row at 318 and push every following row down by one.
$ws.Rows("318:318").Insert()

# Populate the newly inserted row 318 with the new weekly record.
$ws.Cells.Item(318, 1).Value2 = 8
$ws.Cells.Item(318, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(318, 3).Value2 = "Coquimbo"
$ws.Cells.Item(318, 4).Value2 = 45265
$ws.Cells.Item(318, 5).Value2 = 4
$ws.Cells.Item(318, 6).Value2 = 100112001
$ws.Cells.Item(318, 7).Value2 = "Berenjena"
$ws.Cells.Item(318, 8).Value2 = "Sin especificar"
$ws.Cells.Item(318, 9).Value2 = "Primera"
$ws.Cells.Item(318, 10).Value2 = 500
$ws.Cells.Item(318, 11).Value2 = 11000
$ws.Cells.Item(318, 12).Value2 = 12000
$ws.Cells.Item(318, 13).Value2 = 11500
$ws.Cells.Item(318, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(318, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(318, 16).Value2 = 230
$ws.Cells.Item(318, 17).Value2 = 50
$ws.Cells.Item(318, 18).Value2 = "Hortaliza"

# Make sure the date column keeps the same datetime number format as the
# rest of column D (Insert already copies row 317's formatting down, but
# set it explicitly too so it's not dependent on that behaviour).
$ws.Cells.Item(318, 4).NumberFormat = $ws.Cells.Item(317, 4).NumberFormat
